# "Added Indian MF 1st Stab"
#
# This MarketBeat-rank sheet has one row per analyst/firm (col A) and one
# column per weekly report date (row 1, newest date first). This edit adds
# 9 more weekly report dates' worth of data:
#   - Row 1 gets 9 new date labels inserted at the front (cols B..J); the
#     dates that used to start at column B slide right, keeping their text,
#     to make room.
#   - Every other used row just grows by 9 columns, appended right after
#     that row's last populated cell, repeating that row's existing rating
#     value (e.g. "UN") - matching how each analyst row already trails off
#     early if that firm doesn't have history going back as far as others.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDates = @("Sep_08", "Aug_25", "Aug_04", "Jul_23", "Jul_17", "Jul_07", "Jun_30", "Jun_24", "Jun_16")
$shiftBy = $newDates.Count

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row()
$lastRow = $firstRow + $usedRange.Rows.Count() - 1
$firstCol = $usedRange.Column()
$lastColOverall = $firstCol + $usedRange.Columns.Count() - 1

# --- Row 1 (header / date row): shift the existing dates right by 9 columns,
# working from the rightmost column down so we never clobber a value before
# it has been copied, then drop the new dates into the vacated front cols --
$headerRow = $firstRow
$oldHeaderLastCol = $ws.Cells.Item($headerRow, $lastColOverall + 1000).End(-4159).Column()

for ($c = $oldHeaderLastCol; $c -ge $firstCol + 1; $c--) {
    $val = $ws.Cells.Item($headerRow, $c).Value()
    $ws.Cells.Item($headerRow, $c + $shiftBy).Value = $val
}
for ($i = 0; $i -lt $newDates.Count; $i++) {
    $ws.Cells.Item($headerRow, $firstCol + 1 + $i).Value = $newDates[$i]
}

# --- Every other row: append 9 cells after that row's own last populated
# column, repeating whatever value already sits in that last cell ----------
for ($r = $headerRow + 1; $r -le $lastRow; $r++) {
    $rowLastCol = $ws.Cells.Item($r, $lastColOverall + 1000).End(-4159).Column()
    if ($rowLastCol -ge $firstCol) {
        $fillValue = $ws.Cells.Item($r, $rowLastCol).Value()
        for ($i = 1; $i -le $shiftBy; $i++) {
            $ws.Cells.Item($r, $rowLastCol + $i).Value = $fillValue
        }
    }
}
